$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3-8 get reshuffled (same records, different row order).
# New row N gets the D/L/M/N/O/P/S values that used to sit in a different row.
# Mapping (new row -> source row in the original layout):
#   3 <- 4
#   4 <- 7
#   5 <- 8
#   6 <- 3
#   7 <- 5
#   8 <- 6

$rows = @(
    @{ D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; S = 1861 },
    @{ D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; S = 1167 },
    @{ D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; S = 1000 },
    @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; S = 1028 },
    @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; S = 806  },
    @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; S = 667  }
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 4).Value = $data.D        # D - Fecha
    $ws.Cells.Item($r, 12).Value = $data.L       # L - Calidad
    $ws.Cells.Item($r, 13).Value = $data.M       # M - Volumen
    $ws.Cells.Item($r, 14).Value = $data.N       # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.O       # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $data.P       # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $data.S       # S - Precio $/Kg
}
